$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversion del dia" text in A1 with the new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.51 = 39971.48 pesos`n✅ 39971.48 pesos = 9.47 = 966.14 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the rate figures in columns N and O ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 105.2
$ws2.Range("O10").Value = 4205
$ws2.Range("N12").Value = 4219.99
$ws2.Range("O12").Value = 102
